$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The original row 2 (appid=com.hamxa.shaynachim / keyword=bitcoin guide /
# email=vikicrestina@gmail.com / recovery=cristianjohn1222@gmail.com review)
# is being removed entirely; every row below it shifts up by one.

# Hyperlinks don't auto-adjust with a row delete in this host, so clear them
# all first and rebuild them at their new (post-delete) locations below.
$ws.Hyperlinks.Delete()

$ws.Rows("2:2").Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:bittonnir12@gmail.com", "", "", "bittonnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:nevilgreen12@gmail.com", "", "", "nevilgreen12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:shmulmaor2@gmail.com", "", "", "shmulmaor2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:rocketaso@gmail.com", "", "", "rocketaso@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:armonravid@gmail.com", "", "", "armonravid@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:ronoren61@gmail.com", "", "", "ronoren61@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:nitanoren23@gmail.com", "", "", "nitanoren23@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:danfogel100@gmail.com", "", "", "danfogel100@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:avishaybar12@gmail.com", "", "", "avishaybar12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:danfogel100@gmail.com", "", "", "danfogel100@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:avishaybar12@gmail.com", "", "", "avishaybar12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D13"), "mailto:jorjkluni03@gmail.com", "", "", "jorjkluni03@gmail.com")
